$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 117
$ws1.Range("F3").Value = 11
$ws1.Range("F4").Value = 77

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 117
$ws4.Range("F3").Value = 11
$ws4.Range("F4").Value = 77
